$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Header row 1: a new "HSN/SAC" column is inserted before the old
# column F, shifting the old F..J header labels one column to the
# right (into G..K). The old K1 ("HSN/SAC") duplicate is dropped.
# ------------------------------------------------------------------
$ws.Range("F1").Value = "HSN/SAC"
$ws.Range("G1").Value = "Company Name"
$ws.Range("H1").Value = "Invoice No"
$ws.Range("I1").Value = "Date of Invoice"
$ws.Range("J1").Value = "GSTIN NO"
$ws.Range("K1").Value = "GSTIN"
# L1 ("Shipped to") and onward are untouched.

# ------------------------------------------------------------------
# Row 2: Goods Description trimmed, Rate/Amount corrected, and the
# same column shift (F..J -> G..K) applied to the data; the old
# second GSTIN (J2) is dropped, the HSN/SAC code moves into F2, and
# everything from K2 onward is cleared out.
# ------------------------------------------------------------------
$ws.Range("A2").Value = "RICE"

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30"

$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "441"

$ws.Range("F2").NumberFormat = "@"
$ws.Range("F2").Value = "10063090"

$ws.Range("G2").Value = "TANISHQ AGRO INDUSTRIES"

$ws.Range("H2").NumberFormat = "@"
$ws.Range("H2").Value = "595"

$ws.Range("I2").NumberFormat = "@"
$ws.Range("I2").Value = "01-08-2023"

$ws.Range("J2").Value = "27AAKPW5971G1Z1"

$ws.Range("K2:P2").ClearContents()

# ------------------------------------------------------------------
# Row 3: identical treatment to row 2.
# ------------------------------------------------------------------
$ws.Range("A3").Value = "RICE"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "10"

$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "441"

$ws.Range("F3").NumberFormat = "@"
$ws.Range("F3").Value = "10063090"

$ws.Range("G3").Value = "TANISHQ AGRO INDUSTRIES"

$ws.Range("H3").NumberFormat = "@"
$ws.Range("H3").Value = "595"

$ws.Range("I3").NumberFormat = "@"
$ws.Range("I3").Value = "01-08-2023"

$ws.Range("J3").Value = "27AAKPW5971G1Z1"

$ws.Range("K3:P3").ClearContents()
